# Commit: "Added Indian MF 1st Stab"
#
# This workbook tracks analyst-rating snapshots over time. Each week a new
# "as of" date column is inserted at the front of the table (column B),
# pushing older date columns to the right, and a new filler value ("UN")
# is appended at the tail of every existing data row for that new period.
#
# In this edit, 9 new weekly snapshots are added:
#   Jun_16, Jun_24, Jun_30, Jul_07, Jul_17, Jul_23, Aug_04, Aug_25, Sep_08
# (oldest-of-the-new-batch first), which end up placed newest-first in
# columns B..J of row 1, with the previously-existing headers (old B1:V1)
# shifted right by 9 columns (now K1:AE1). Existing data rows keep their
# original cells untouched and simply grow by 9 more "UN" cells appended
# after their current last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shiftCount = 9

# ---------------------------------------------------------------------
# Step 1: shift the existing header row (row 1) to the right by 9 cols.
# Row 1 currently has data in columns B (2) .. V (22).
# We copy from the rightmost column down to the leftmost so that we never
# overwrite a source cell before it has been read.
# ---------------------------------------------------------------------
$firstCol = 2   # column B
$lastCol  = 22  # column V

for ($c = $lastCol; $c -ge $firstCol; $c--) {
    $srcVal = $ws.Cells.Item(1, $c).Value()
    $ws.Cells.Item(1, $c + $shiftCount).Value = $srcVal
}

# ---------------------------------------------------------------------
# Step 2: write the 9 new date headers into the freshly vacated
# columns B..J (2..10), newest date first.
# ---------------------------------------------------------------------
$newDates = @("Sep_08", "Aug_25", "Aug_04", "Jul_23", "Jul_17", "Jul_07", "Jun_30", "Jun_24", "Jun_16")

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $ws.Cells.Item(1, $firstCol + $i).Value = $newDates[$i]
}

# ---------------------------------------------------------------------
# Step 3: for every data row, append 9 more filler cells ("UN") right
# after whatever its current last populated column is. Row 1 is skipped
# since it was already handled above. Rows with no data are skipped too.
# ---------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $lastColInRow = $ws.Cells.Item($r, $ws.Columns.Count).End(-4159).Column
    if ($lastColInRow -lt $firstCol) {
        continue
    }
    for ($i = 1; $i -le $shiftCount; $i++) {
        $ws.Cells.Item($r, $lastColInRow + $i).Value = "UN"
    }
}

# ---------------------------------------------------------------------
# Step 4: extend the custom column width formatting that already existed
# on columns C..V so the 9 newly used columns (W..AE) match it too.
# ---------------------------------------------------------------------
$refWidth = $ws.Columns.Item(3).ColumnWidth
for ($c = 23; $c -le 31; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $refWidth
}
